$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric cell value updates
$ws.Range("F122").Value2 = 24
$ws.Range("G122").Value2 = 1702.08
$ws.Range("B147").Value2 = 106524.9
$ws.Range("B151").Value2 = 65258
$ws.Range("F151").Value2 = 2
$ws.Range("G151").Value2 = 64287.16
$ws.Range("B152").Value2 = 64196
$ws.Range("F152").Value2 = 1
$ws.Range("G152").Value2 = 32143.58
$ws.Range("F182").Value2 = 15
$ws.Range("G182").Value2 = 985.05
$ws.Range("B184").Value2 = 33619.23
$ws.Range("F227").Value2 = 6
$ws.Range("G227").Value2 = 276.12
$ws.Range("B228").Value2 = 57001
$ws.Range("E228").Value2 = 103.53
$ws.Range("F228").Value2 = 3
$ws.Range("G228").Value2 = 259.95
$ws.Range("B229").Value2 = 64693
$ws.Range("E229").Value2 = 92.13
$ws.Range("F229").Value2 = 1
$ws.Range("G229").Value2 = 86.65000000000001
$ws.Range("F234").Value2 = 36
$ws.Range("G234").Value2 = 1826.64
$ws.Range("B237").Value2 = 57004
$ws.Range("F237").Value2 = 5
$ws.Range("G237").Value2 = 410
$ws.Range("B238").Value2 = 63255
$ws.Range("F238").Value2 = 105
$ws.Range("G238").Value2 = 8610
$ws.Range("B251").Value2 = 110099.78
$ws.Range("F268").Value2 = 4
$ws.Range("G268").Value2 = 185.96
$ws.Range("B274").Value2 = 11976.37
$ws.Range("F333").Value2 = 103
$ws.Range("G333").Value2 = 8740.58
$ws.Range("F342").Value2 = 50
$ws.Range("G342").Value2 = 6360.5
$ws.Range("B371").Value2 = 66194
$ws.Range("F371").Value2 = 46
$ws.Range("G371").Value2 = 3941.28
$ws.Range("B372").Value2 = 64983
$ws.Range("F372").Value2 = 6
$ws.Range("G372").Value2 = 514.08
$ws.Range("B376").Value2 = 199089.81
$ws.Range("B379").Value2 = 63565
$ws.Range("E379").Value2 = 109.19
$ws.Range("F379").Value2 = 60
$ws.Range("G379").Value2 = 6162.6
$ws.Range("B380").Value2 = 61610
$ws.Range("E380").Value2 = 122.71
$ws.Range("F380").Value2 = -58
$ws.Range("G380").Value2 = -5957.18
$ws.Range("B401").Value2 = 60325
$ws.Range("E401").Value2 = 151.57
$ws.Range("F401").Value2 = -102
$ws.Range("G401").Value2 = -12939.72
$ws.Range("B402").Value2 = 63560
$ws.Range("E402").Value2 = 134.87
$ws.Range("F402").Value2 = 1
$ws.Range("G402").Value2 = 126.86
$ws.Range("F457").Value2 = 80
$ws.Range("G457").Value2 = 18417.6
$ws.Range("F458").Value2 = 72
$ws.Range("G458").Value2 = 17287.2
$ws.Range("B462").Value2 = 139265.4
$ws.Range("F468").Value2 = 186
$ws.Range("G468").Value2 = 11290.2
$ws.Range("B474").Value2 = 109983.53
$ws.Range("B484").Value2 = 58047
$ws.Range("D484").Value2 = 105.54
$ws.Range("E484").Value2 = 126.1
$ws.Range("F484").Value2 = 35
$ws.Range("G484").Value2 = 3693.9
$ws.Range("B485").Value2 = 47097
$ws.Range("D485").Value2 = 112.28
$ws.Range("E485").Value2 = 134.16
$ws.Range("F485").Value2 = 15
$ws.Range("G485").Value2 = 1684.2
$ws.Range("F506").Value2 = 121
$ws.Range("G506").Value2 = 4151.51
$ws.Range("B523").Value2 = 212616.76
$ws.Range("B563").Value2 = 45718
$ws.Range("E563").Value2 = 19.38
$ws.Range("F563").Value2 = -294
$ws.Range("G563").Value2 = -4768.68
$ws.Range("B564").Value2 = 64927
$ws.Range("E564").Value2 = 17.26
$ws.Range("F564").Value2 = 106
$ws.Range("G564").Value2 = 1719.32
$ws.Range("B570").Value2 = 45702
$ws.Range("E570").Value2 = 31.43
$ws.Range("F570").Value2 = -215
$ws.Range("G570").Value2 = -5654.5
$ws.Range("B571").Value2 = 64919
$ws.Range("E571").Value2 = 27.97
$ws.Range("F571").Value2 = 61
$ws.Range("G571").Value2 = 1604.3
$ws.Range("B573").Value2 = 65067
$ws.Range("E573").Value2 = 15.65
$ws.Range("F573").Value2 = 126
$ws.Range("G573").Value2 = 1855.98
$ws.Range("B574").Value2 = 53595
$ws.Range("E574").Value2 = 17.61
$ws.Range("F574").Value2 = -335
$ws.Range("G574").Value2 = -4934.55
$ws.Range("F585").Value2 = 42
$ws.Range("G585").Value2 = 2925.3
$ws.Range("B588").Value2 = 49945.92
$ws.Range("F612").Value2 = 159
$ws.Range("G612").Value2 = 44990.64
$ws.Range("B616").Value2 = 163813.97
$ws.Range("B644").Value2 = 53319
$ws.Range("E644").Value2 = 310.64
$ws.Range("F644").Value2 = -6
$ws.Range("G644").Value2 = -1643.52
$ws.Range("B645").Value2 = 64810
$ws.Range("E645").Value2 = 291.22
$ws.Range("F645").Value2 = 2
$ws.Range("G645").Value2 = 547.84
$ws.Range("F692").Value2 = 160
$ws.Range("G692").Value2 = 4262.4
$ws.Range("B696").Value2 = 216721.79
$ws.Range("B834").Value2 = 65362
$ws.Range("F834").Value2 = 2
$ws.Range("G834").Value2 = 81.73999999999999
$ws.Range("B835").Value2 = 65079
$ws.Range("F835").Value2 = 6
$ws.Range("G835").Value2 = 245.22
$ws.Range("B948").Value2 = 5989701.5
$ws.Range("B949").Value2 = 5989701.5

# Text cell value updates (item name casing fix due to row data swap)
$ws.Range("C371").Value2 = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("C372").Value2 = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
